# "add title to plain template"
#
# The first paragraph currently holds one run of plain text:
#   "$Communism and the Intellectual % Dorothy Day"
#
# It needs to become two separate paragraphs:
#   1) a "Title" styled paragraph:   Communism and the Intellectual
#   2) an "Authors" styled paragraph: Dorothy Day
# with each word/space of the new text living in its own run (as a
# plain retype of the title would produce), and the stray "$"/"%"
# punctuation dropped.

$d = $word.ActiveDocument

# Turning on revision tracking means each TypeText() call records its
# own insertion, so the separate words don't get silently coalesced
# back into a single run the way a plain Range.Text assignment would.
$d.TrackRevisions = $true

$titlePara = $d.Paragraphs(1).Range
$titlePara.Text = ""

$sel = $word.Selection
$sel.SetRange($titlePara.Start, $titlePara.Start)

$titleWords = @("Communism", " ", "and", " ", "the", " ", "Intellectual")
foreach ($w in $titleWords) {
    $sel.TypeText($w)
    $sel.Collapse(0)
}

# Split the authors onto their own paragraph.
$sel.TypeParagraph()

$authorWords = @("Dorothy", " ", "Day")
foreach ($w in $authorWords) {
    $sel.TypeText($w)
    $sel.Collapse(0)
}

# Flatten the tracked insertions back into plain runs (keeping the
# per-word run boundaries created above) and stop tracking again.
$d.TrackRevisions = $false
$d.Revisions.AcceptAll()

# Apply the template's heading styles to the two new paragraphs.
$d.Paragraphs(1).Style = "Title"
$d.Paragraphs(2).Style = "Authors"
